# Applies the "examples/biochemical_models/data_copy.xlsx" restructuring:
#  - delete the "!_Table of contents" sheet
#  - rename table-metadata keys (TableType->Type, ModelId->Id, ModelName->Name)
#    and bump the embedded generation Date
#  - shorten attribute-type names in the "!_Schema" sheet
#    (XxxAttribute -> Xxx)
#  - fix up the "Model:1" -> "!Model:1" wording in comments + data validations
#    on "!Compound" and "!Reaction"
#  - make "!_Schema" the new first/active sheet

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the table-of-contents sheet entirely.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("!_Table of contents").Delete()

# ---------------------------------------------------------------------------
# 2. "!_Schema" sheet: metadata title + attribute type names.
# ---------------------------------------------------------------------------
$wsSchema = $wb.Worksheets.Item("!_Schema")

$wsSchema.Range("A1").Value = "!!ObjTables Type='Schema' Description='Table/model and column/attribute definitions' Date='2019-10-10 02:11:40' ObjTablesVersion='0.0.8'"

$attrCells = @("D4","D5","D6","D7","D8","D10","D11","D13","D14","D15","D16","D17","D18","D19")
foreach ($addr in $attrCells) {
    $rng = $wsSchema.Range($addr)
    $old = $rng.Value()
    $new = $old
    if ($old -eq "SlugAttribute") { $new = "Slug" }
    elseif ($old -eq "StringAttribute") { $new = "String" }
    elseif ($old -eq "BooleanAttribute") { $new = "Boolean" }
    elseif ($old -eq "ManyToOneAttribute") { $new = "ManyToOne" }
    if ($new -ne $old) {
        $rng.Value = $new
    }
}

# Make it the first / active sheet (matches tabSelected="1" in the diff).
$wsSchema.Move($wb.Worksheets.Item(1))
$wsSchema.Activate()

# ---------------------------------------------------------------------------
# 3. "!Compound" sheet: metadata title, comment + data validation wording.
# ---------------------------------------------------------------------------
$wsCompound = $wb.Worksheets.Item("!Compound")
$wsCompound.Range("A1").Value = "!!ObjTables Type='Data' Id='Compound' Description='Compound' Name='Compound' Date='2019-10-10 02:11:40' ObjTablesVersion='0.0.8'"

$cCompound = $wsCompound.Range("A2").Comment
$cCompound.Text("Select a value from ""!Model:1"" or blank.")

$dvCompound = $wsCompound.Range("A3:A7").Validation
$dvCompound.Modify(3, 2, 1, "'!Model'!`$B`$1:`$XFD`$1")
$dvCompound.ErrorTitle = "Model"
$dvCompound.ErrorMessage = "Value must be a value from ""!Model:1"" or blank."
$dvCompound.InputTitle = "Model"
$dvCompound.InputMessage = "Select a value from ""!Model:1"" or blank."
$dvCompound.IgnoreBlank = $true
$dvCompound.ShowInput = $true
$dvCompound.ShowError = $true

# ---------------------------------------------------------------------------
# 4. "!Model" sheet: metadata title only.
# ---------------------------------------------------------------------------
$wsModel = $wb.Worksheets.Item("!Model")
$wsModel.Range("A1").Value = "!!ObjTables Type='Data' Id='Model' Description='Model' Name='Model' Date='2019-10-10 02:11:40' ObjTablesVersion='0.0.8'"

# ---------------------------------------------------------------------------
# 5. "!Reaction" sheet: metadata title, comment + data validation wording.
# ---------------------------------------------------------------------------
$wsReaction = $wb.Worksheets.Item("!Reaction")
$wsReaction.Range("A1").Value = "!!ObjTables Type='Data' Id='Reaction' Description='Reaction' Name='Reaction' Date='2019-10-10 02:11:40' ObjTablesVersion='0.0.8'"

$cReaction = $wsReaction.Range("A2").Comment
$cReaction.Text("Select a value from ""!Model:1"" or blank.")

$dvReaction = $wsReaction.Range("A3:A4").Validation
$dvReaction.Modify(3, 2, 1, "'!Model'!`$B`$1:`$XFD`$1")
$dvReaction.ErrorTitle = "Model"
$dvReaction.ErrorMessage = "Value must be a value from ""!Model:1"" or blank."
$dvReaction.InputTitle = "Model"
$dvReaction.InputMessage = "Select a value from ""!Model:1"" or blank."
$dvReaction.IgnoreBlank = $true
$dvReaction.ShowInput = $true
$dvReaction.ShowError = $true
